$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 contains the simulation-tool column labels in B1:I1.
# Rename each of them by appending "(APH)" to the existing label,
# e.g. "ESP" -> "ESP(APH)", "BLAST" -> "BLAST(APH)", etc.
$headerCells = @("B1", "C1", "D1", "E1", "F1", "G1", "H1", "I1")

foreach ($cellAddr in $headerCells) {
    $cell = $ws.Range($cellAddr)
    $currentText = $cell.Value2
    $cell.Value = "$currentText(APH)"
}
